$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the same date/time number format used by existing rows (A2:A3) to the new date cells
$ws.Range("A4:A17").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A4").Value = 45211
$ws.Range("B4").Value = "Victory Drones"
$ws.Range("C4").Value = 350000

$ws.Range("A5").Value = 45211
$ws.Range("B5").Value = "General donations"
$ws.Range("C5").Value = 3000

$ws.Range("A6").Value = 45194
$ws.Range("B6").Value = "General donations"
$ws.Range("C6").Value = 15000

$ws.Range("A7").Value = 45192
$ws.Range("B7").Value = "General donations"
$ws.Range("C7").Value = 3000

$ws.Range("A8").Value = 45191
$ws.Range("B8").Value = "General donations"
$ws.Range("C8").Value = 3089.78

$ws.Range("A9").Value = 45191
$ws.Range("B9").Value = "General donations"
$ws.Range("C9").Value = 12000

$ws.Range("A10").Value = 45189
$ws.Range("B10").Value = "General donations"
$ws.Range("C10").Value = 3000

$ws.Range("A11").Value = 45163
$ws.Range("B11").Value = "General donations"
$ws.Range("C11").Value = 10000

$ws.Range("A12").Value = 45161
$ws.Range("B12").Value = "General donations"
$ws.Range("C12").Value = 3500

$ws.Range("A13").Value = 45160
$ws.Range("B13").Value = "General donations"
$ws.Range("C13").Value = 10000

$ws.Range("A14").Value = 45152
$ws.Range("B14").Value = "General donations"
$ws.Range("C14").Value = 250000

$ws.Range("A15").Value = 45148
$ws.Range("B15").Value = "General donations"
$ws.Range("C15").Value = 840000

$ws.Range("A16").Value = 45098
$ws.Range("B16").Value = "General donations"
$ws.Range("C16").Value = 3456

$ws.Range("A17").Value = 45068
$ws.Range("B17").Value = "General donations"
$ws.Range("C17").Value = 3673
